# Aggiornato file Ore con orario del 18/12 e del 19/12.
# Manca ancora la rendicontazione del 19/12 per Mirko.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New activity log rows -------------------------------------------------
# Row 17: 18/12/2017 - Giovanni - "Scrittura metodo combina, crea, ricicla" - 3.5h
$ws.Range("A17").Value = 43087
$ws.Range("B17").Value = "Giovanni"
$ws.Range("C17").Value = "Scrittura metodo combina, crea, ricicla"
$ws.Range("D17").Value = 0.14583333333333334
$ws.Rows.Item(17).RowHeight = 29.15

# Row 18: 18/12/2017 - Mirko - "Scrittura metodo combina, crea, ricicla" - 3.5h
$ws.Range("A18").Value = 43087
$ws.Range("B18").Value = "Mirko"
$ws.Range("C18").Value = "Scrittura metodo combina, crea, ricicla"
$ws.Range("D18").Value = 0.14583333333333334
$ws.Rows.Item(18).RowHeight = 29.15

# Row 19: 19/12/2017 - Giovanni - "Debug di crea :(" - 2.5h (Mirko not yet logged for this day)
$ws.Range("A19").Value = 43088
$ws.Range("B19").Value = "Giovanni"
$ws.Range("C19").Value = "Debug di crea :("
$ws.Range("D19").Value = 0.10416666666666667

# --- Totals now roll up in hours that can exceed 24, switch to elapsed-time format
$ws.Range("F2").NumberFormat = "[h]:mm:ss"
$ws.Range("G2").NumberFormat = "[h]:mm:ss"

# --- View: scroll back to the top and land the selection on G2 -------------
$ws.Range("G2").Select()
